$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new testing statistics values for row 3
$ws.Range("A3").Value = 1000
$ws.Range("B3").Value = 44286
$ws.Range("C3").Value = 11
$ws.Range("D3").Value = 16
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 0.968

# Match formatting of the row above (row 2)
$ws.Range("A2:F2").Copy()
$ws.Range("A3:F3").PasteSpecial(-4122)

# Row 2 uses a taller row height (larger font) - match it for the new row
$ws.Rows.Item(3).RowHeight = 15.6

$ws.Range("E14").Select()
